# Appends new sensor-log rows to the "PIR" sheet (rows 101-126) and the
# "Humidity" sheet (rows 56-76), matching the latest export from the
# logging pipeline.

$wb = $excel.ActiveWorkbook

function Set-LogRow {
    param(
        $ws,
        [int]$row,
        [string]$date,
        [string]$timestamp,
        [string]$hour,
        [string]$location,
        [string]$value,
        [string]$status
    )

    # Column A holds a date-shaped string ("2026-02-01"). Excel's COM value
    # setter auto-detects that shape and silently coerces it to a date
    # serial number, so the cell is forced to Text format first to keep the
    # literal string.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $date

    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location

    # Column E sometimes holds a percentage-shaped string (e.g. "79.4%"),
    # which Excel would otherwise auto-convert to a numeric percentage.
    $cellE = $ws.Cells.Item($row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $value

    $ws.Cells.Item($row, 6).Value = $status
}

# --- PIR sheet: rows 101-126 -------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")

$pirRows = @(
    @("2026-02-01","13:58:26","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:26","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:27","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:27","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:32","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:32","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:37","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:38","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:42","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:42","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:47","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:47","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:52","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:52","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:57","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:58:57","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:59:02","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:59:02","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:59:07","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:59:07","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:59:12","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:59:12","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:59:17","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:59:18","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:59:22","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:59:23","13:00","Bathroom","No Motion","Inactive")
)

$row = 101
foreach ($r in $pirRows) {
    Set-LogRow $wsPir $row $r[0] $r[1] $r[2] $r[3] $r[4] $r[5]
    $row++
}

# --- Humidity sheet: rows 56-76 -----------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")

$humidityRows = @(
    @("2026-02-01","13:58:26","13:00","Bathroom","79.4%","Active"),
    @("2026-02-01","13:58:26","13:00","Bathroom","78.5%","Active"),
    @("2026-02-01","13:58:27","13:00","Bathroom","79.4%","Active"),
    @("2026-02-01","13:58:29","13:00","Bathroom","79.4%","Active"),
    @("2026-02-01","13:58:34","13:00","Bathroom","78.5%","Active"),
    @("2026-02-01","13:58:36","13:00","Bathroom","79.5%","Active"),
    @("2026-02-01","13:58:39","13:00","Bathroom","79.4%","Active"),
    @("2026-02-01","13:58:41","13:00","Bathroom","79.5%","Active"),
    @("2026-02-01","13:58:46","13:00","Bathroom","79.5%","Active"),
    @("2026-02-01","13:58:49","13:00","Bathroom","79.4%","Active"),
    @("2026-02-01","13:58:51","13:00","Bathroom","79.4%","Active"),
    @("2026-02-01","13:58:55","13:00","Bathroom","78.5%","Active"),
    @("2026-02-01","13:58:56","13:00","Bathroom","79.4%","Active"),
    @("2026-02-01","13:58:59","13:00","Bathroom","79.4%","Active"),
    @("2026-02-01","13:59:04","13:00","Bathroom","78.4%","Active"),
    @("2026-02-01","13:59:10","13:00","Bathroom","79.4%","Active"),
    @("2026-02-01","13:59:11","13:00","Bathroom","79.3%","Active"),
    @("2026-02-01","13:59:15","13:00","Bathroom","78.3%","Active"),
    @("2026-02-01","13:59:20","13:00","Bathroom","77.8%","Active"),
    @("2026-02-01","13:59:21","13:00","Bathroom","79.4%","Active"),
    @("2026-02-01","13:59:25","13:00","Bathroom","78.3%","Active")
)

$row = 56
foreach ($r in $humidityRows) {
    Set-LogRow $wsHumidity $row $r[0] $r[1] $r[2] $r[3] $r[4] $r[5]
    $row++
}
